# Mejora fase de elaboracion
#
# Word's proofing pass marks "Mg.Sc." and "cc." with proofErr
# (spellStart/spellEnd) bookmarks because the speller treats the
# abbreviation before the final period as a separate "word". This
# splits each of those two runs into two runs - the abbreviation
# (wrapped in proofErr spellStart/spellEnd) and the trailing text -
# using Range.InsertXML so the proofErr markers land exactly where
# Word would place them.

$d = $word.ActiveDocument

$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'
$pkgNs = 'http://schemas.microsoft.com/office/2006/xmlPackage'

function Split-RunWithProofErr($FullText, $FirstPart, $SecondPart, $RunPropsXml) {
    $d2 = $word.ActiveDocument
    $target = $d2.Content
    $target.Find.ClearFormatting()
    $found = $target.Find.Execute($FullText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found -or -not $target.Find.Found) {
        throw "Could not find text '$FullText'"
    }

    $start = $target.Start
    $firstRange = $d2.Range($start, $start + $FirstPart.Length)
    $secondRange = $d2.Range($start + $FirstPart.Length, $start + $FullText.Length)

    if ($firstRange.Text -ne $FirstPart) {
        throw "First-part mismatch: got [$($firstRange.Text)] expected [$FirstPart]"
    }
    if ($secondRange.Text -ne $SecondPart) {
        throw "Second-part mismatch: got [$($secondRange.Text)] expected [$SecondPart]"
    }

    $secondSpace = ''
    if ($SecondPart.Length -gt 0 -and ($SecondPart[0] -eq ' ' -or $SecondPart[$SecondPart.Length - 1] -eq ' ')) {
        $secondSpace = ' xml:space="preserve"'
    }

    $fragment = '<w:proofErr w:type="spellStart"/>' +
                '<w:r>' + $RunPropsXml + '<w:t>' + $FirstPart + '</w:t></w:r>' +
                '<w:proofErr w:type="spellEnd"/>' +
                '<w:r>' + $RunPropsXml + '<w:t' + $secondSpace + '>' + $SecondPart + '</w:t></w:r>'

    $payload = '<pkg:package xmlns:pkg="' + $pkgNs + '">' +
                 '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
                   '<pkg:xmlData>' +
                     '<w:p xmlns:w="' + $wNs + '">' + $fragment + '</w:p>' +
                   '</pkg:xmlData>' +
                 '</pkg:part>' +
               '</pkg:package>'

    # Re-fetch the (possibly re-seated) range right before the edit and
    # replace exactly the full matched span with the new run structure.
    $editRange = $d2.Range($start, $start + $FullText.Length)
    $editRange.InsertXML($payload)
}

$runProps = '<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-EC"/></w:rPr>'

Split-RunWithProofErr "Mg.Sc." "Mg.Sc" "." $runProps
Split-RunWithProofErr "cc. Archivo, Elisa Orellana" "cc." " Archivo, Elisa Orellana" $runProps

Write-Output "done"
